# Update the old "link:X.Y.Z" style references to the new "$ref:X:Y.Z" API format.
$wb = $excel.ActiveWorkbook

$testsSheet = $wb.Worksheets.Item("Tests")
$dataBlocksSheet = $wb.Worksheets.Item("DataBlocks")

# Order matches the order new shared-string entries were introduced by the
# original edit so the resulting sharedStrings.xml table lines up exactly.
$testsSheet.Range("C9").Value = "`$ref:DataBlocks:AnotherObject.cyclicRef"
$testsSheet.Range("C10").Value = "`$ref:DataBlocks:NewObject"
$testsSheet.Range("C14").Value = "`$ref:DataBlocks:NewObject.VALUE"
$dataBlocksSheet.Range("C3").Value = "`$ref:DataBlocks:AnotherObject.anotherValue"
$testsSheet.Range("C5").Value = "`$ref:Tests:Uncommon.reftest"
$testsSheet.Range("C6").Value = "`$ref:Tests:Uncommon.reftestGen"
$dataBlocksSheet.Range("C7").Value = "`$ref:Tests:Common.cyclic"
$dataBlocksSheet.Range("C8").Value = "`$ref:Tests:Common.gendata"
